$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 89
$ws.Range("D7").Value = 84
$ws.Range("D8").Value = 79

$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 14.25

$ws.Range("H2:H13").Formula = "=SUM(A2:G2)"

$ws.Range("D9").Select()
